$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.150.11"
$ws.Range("E2").Value = "  -2.14%  "

# Row 3
$ws.Range("D3").Value = "1.577.40"

# Row 4
$ws.Range("E4").Value = "  -0.38%  "

# Row 5
$ws.Range("D5").Value = "'209.21"
$ws.Range("E5").Value = "  -1.12%  "

# Row 6
$ws.Range("D6").Value = "'0.497"
$ws.Range("E6").Value = "  -3.12%  "

# Row 7
$ws.Range("E7").Value = "  -0.33%  "

# Row 8
$ws.Range("E8").Value = "  -0.56%  "

# Row 9
$ws.Range("E9").Value = "  -1.47%  "

# Row 10
$ws.Range("D10").Value = "'19.54"
$ws.Range("E10").Value = "  -0.59%  "

# Row 11
$ws.Range("D11").Value = "'0.0844"

# Row 12
$ws.Range("D12").Value = "1.798.96"
$ws.Range("E12").Value = "  -1.45%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.583.24"
$ws.Range("E13").Value = "  -1.04%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'4.05"
$ws.Range("E14").Value = "  -0.05%  "

# Row 15
$ws.Range("E15").Value = "  -1.85%  "

# Row 16
$ws.Range("D16").Value = "'64.44"
$ws.Range("E16").Value = "  -0.94%  "

# Row 17
$ws.Range("D17").Value = "26.150.04"
$ws.Range("E17").Value = "  -2.05%  "

# Row 18
$ws.Range("E18").Value = "  -1.11%  "

# Row 19
$ws.Range("E19").Value = "  +1.47%  "

# Row 20
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "'1.00"
$ws.Range("E20").Value = "  -0.34%  "

# Row 21
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'207.41"
$ws.Range("E21").Value = "  -1.29%  "

# Row 22
$ws.Range("E22").Value = "  -0.94%  "

# Row 23
$ws.Range("E23").Value = "  -3.05%  "

# Row 24
$ws.Range("D24").Value = "'8.89"
$ws.Range("E24").Value = "  -1.11%  "

# Row 25
$ws.Range("D25").Value = "'143.90"
$ws.Range("E25").Value = "  +0.21%  "

# Row 26
$ws.Range("E26").Value = "  -0.36%  "

# Row 27
$ws.Range("E27").Value = "  -1.35%  "

# Row 28
$ws.Range("E28").Value = "  -1.64%  "

# Row 29
$ws.Range("E29").Value = "  -0.78%  "

# Row 30
$ws.Range("D30").Value = "'0.0506"
$ws.Range("E30").Value = "  -0.65%  "

# Row 31
$ws.Range("E31").Value = "  -1.27%  "

# Row 32
$ws.Range("E32").Value = "  -1.65%  "

# Row 33
$ws.Range("D33").Value = "'2.98"
$ws.Range("E33").Value = "  +0.57%  "

# Row 34
$ws.Range("D34").Value = "1.278.56"
$ws.Range("E34").Value = "  -0.86%  "

# Row 35
$ws.Range("E35").Value = "  -0.70%  "

# Row 36
$ws.Range("D36").Value = "'0.607"
$ws.Range("E36").Value = "  +0.94%  "

# Row 37
$ws.Range("E37").Value = "  -1.03%  "

# Row 38
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "'1.11"
$ws.Range("E38").Value = "  -4.83%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.0165"
$ws.Range("E39").Value = "  -2.73%  "

# Row 40
$ws.Range("E40").Value = "  -1.65%  "

# Row 41
$ws.Range("D41").Value = "'5.55"
$ws.Range("E41").Value = "  +2.93%  "

# Row 42
$ws.Range("E42").Value = "  -2.32%  "

# Row 43
$ws.Range("D43").Value = "'62.44"
$ws.Range("E43").Value = "  -0.74%  "

# Row 44
$ws.Range("E44").Value = "  -2.63%  "

# Row 45
$ws.Range("D45").Value = "1.712.09"
$ws.Range("E45").Value = "  -1.41%  "

# Row 46
$ws.Range("D46").Value = "'88.96"
$ws.Range("E46").Value = "  -1.72%  "

# Row 47
$ws.Range("E47").Value = "  -0.61%  "

# Row 48
$ws.Range("E48").Value = "  -1.77%  "

# Row 49
$ws.Range("E49").Value = "  -1.23%  "

# Row 50
$ws.Range("E50").Value = "  -1.98%  "

# Row 51
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  -0.29%  "
